$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the I column (praclen) values for existing rows 2-5: 4 -> 5
$ws.Range("I2").Value = 5
$ws.Range("I3").Value = 5
$ws.Range("I4").Value = 5
$ws.Range("I5").Value = 5

# Add new row 6 with the new training schedule entry
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 6
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = -5
$ws.Range("G6").Value = -1
$ws.Range("H6").Value = 12
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = "train_dim2_1"

# Update selection to match the final saved state (A5 selected)
$ws.Range("A5").Select()
